$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Title row (row 1): bold, size 14, centered ---
$ws.Range("A1:G1").Font.Bold = $true
$ws.Range("A1:G1").Font.Size = 14
$ws.Range("A1:G1").HorizontalAlignment = -4108
$ws.Rows.Item(1).RowHeight = 18.75

# --- Header row (row 2): bold ---
$ws.Range("A2:G2").Font.Bold = $true

# --- Rewrite data rows 3-18 (new backlog items, Adjust Factor all 0) ---
# Row 15 previously held the bold "Total:" cell; clear that formatting
# before writing the new (non-bold) data into it.
$ws.Range("A15").Font.Bold = $false

$ws.Range("A3").Value = "Welcome Page"
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 0
$ws.Range("D3").Formula = "=B3+(B3*C3)"

$ws.Range("A4").Value = "Login Page"
$ws.Range("B4").Value = 7
$ws.Range("C4").Value = 0
$ws.Range("D4").Formula = "=B4+(B4*C4)"

$ws.Range("A5").Value = "Actual Login Implementation"
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = 0
$ws.Range("D5").Formula = "=B5+(B5*C5)"

$ws.Range("A6").Value = "Register Page"
$ws.Range("B6").Value = 7
$ws.Range("C6").Value = 0
$ws.Range("D6").Formula = "=B6+(B6*C6)"

$ws.Range("A7").Value = "Actual Register Implementation"
$ws.Range("B7").Value = 10
$ws.Range("C7").Value = 0
$ws.Range("D7").Formula = "=B7+(B7*C7)"

$ws.Range("A8").Value = "Admin Account"
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = 0
$ws.Range("D8").Formula = "=B8+(B8*C8)"

$ws.Range("A9").Value = "Create Multiple Accounts"
$ws.Range("B9").Value = 5
$ws.Range("C9").Value = 0
$ws.Range("D9").Formula = "=B9+(B9*C9)"

$ws.Range("A10").Value = "User Profile"
$ws.Range("B10").Value = 3
$ws.Range("C10").Value = 0
$ws.Range("D10").Formula = "=B10+(B10*C10)"

$ws.Range("A11").Value = "Deposits"
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = 0
$ws.Range("D11").Formula = "=B11+(B11*C11)"

$ws.Range("A12").Value = "Withdrawals"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 0
$ws.Range("D12").Formula = "=B12+(B12*C12)"

$ws.Range("A13").Value = "Transaction History"
$ws.Range("B13").Value = 10
$ws.Range("C13").Value = 0
$ws.Range("D13").Formula = "=B13+(B13*C13)"

$ws.Range("A14").Value = "Spending Catergory Report"
$ws.Range("B14").Value = 15
$ws.Range("C14").Value = 0
$ws.Range("D14").Formula = "=B14+(B14*C14)"

$ws.Range("A15").Value = "Income Source Report"
$ws.Range("B15").Value = 15
$ws.Range("C15").Value = 0
$ws.Range("D15").Formula = "=B15+(B15*C15)"

$ws.Range("A16").Value = "Cash Flow Reports"
$ws.Range("B16").Value = 15
$ws.Range("C16").Value = 0
$ws.Range("D16").Formula = "=B16+(B16*C16)"

$ws.Range("A17").Value = "Account Listing Report"
$ws.Range("B17").Value = 15
$ws.Range("C17").Value = 0
$ws.Range("D17").Formula = "=B17+(B17*C17)"

$ws.Range("A18").Value = "Transaction History Report"
$ws.Range("B18").Value = 15
$ws.Range("C18").Value = 0
$ws.Range("D18").Formula = "=B18+(B18*C18)"

# --- Row 19 stays blank (gap between data and total) ---

# --- Total row moved to row 20 ---
$ws.Range("A20").Value = "Total:"
$ws.Range("A20").Font.Bold = $true
$ws.Range("B20").Formula = "=SUM(B3:B18)"
$ws.Range("D20").Formula = "=SUM(D3:D18)"

# --- Selection moves to A15 ---
$ws.Range("A15").Select()
